$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 4501.552
$ws.Range("I86").Value = 1020.12
$ws.Range("J86").Value = 26260.5
$ws.Range("K86").Value = 1020.12
$ws.Range("L86").Value = 26260.5
$ws.Range("M86").Value = 102.88
$ws.Range("N86").Value = -28506.5
$ws.Range("H89").Value = 4501.552
$ws.Range("I89").Value = 1020.12
$ws.Range("J89").Value = 26260.5
$ws.Range("K89").Value = 5100.6
$ws.Range("L89").Value = 131302.5
$ws.Range("M89").Value = 515.3999999999996
$ws.Range("N89").Value = -142534.5
$ws.Range("H116").Value = 16607422
$ws.Range("J116").Value = 5436.5454
$ws.Range("L116").Value = 5436.5454
$ws.Range("N116").Value = -12320.5454
$ws.Range("H131").Value = 1905.8096
$ws.Range("I131").Value = 885.1667
$ws.Range("K131").Value = 2655.5001
$ws.Range("M131").Value = 2384.4999
$ws.Range("H137").Value = 90493.625
$ws.Range("I137").Value = 109534.164
$ws.Range("J137").Value = 2431.125
$ws.Range("K137").Value = 328602.492
$ws.Range("L137").Value = 7293.375
$ws.Range("M137").Value = -326052.492
$ws.Range("N137").Value = -12393.375
$ws.Range("H141").Value = 1808
$ws.Range("I141").Value = 1512.1936
$ws.Range("J141").Value = 3642
$ws.Range("K141").Value = 4536.5808
$ws.Range("L141").Value = 10926
$ws.Range("M141").Value = 643.4192000000003
$ws.Range("N141").Value = -21286

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 9307.329
$ws.Range("I32").Value = 6751.2188
$ws.Range("K32").Value = 6751.2188
$ws.Range("M32").Value = -6464.2188
$ws.Range("H63").Value = 3127491
$ws.Range("I63").Value = 2738.625
$ws.Range("K63").Value = 2738.625
$ws.Range("M63").Value = -2052.625
$ws.Range("H66").Value = 3127491
$ws.Range("I66").Value = 2738.625
$ws.Range("K66").Value = 13693.125
$ws.Range("M66").Value = -10261.125
$ws.Range("H74").Value = 32259444
$ws.Range("I74").Value = 40000668
$ws.Range("J74").Value = 4350
$ws.Range("K74").Value = 40000668
$ws.Range("L74").Value = 4350
$ws.Range("M74").Value = -39999794
$ws.Range("N74").Value = -6098
$ws.Range("H77").Value = 32259444
$ws.Range("I77").Value = 40000668
$ws.Range("J77").Value = 4350
$ws.Range("K77").Value = 200003340
$ws.Range("L77").Value = 21750
$ws.Range("M77").Value = -199998972
$ws.Range("N77").Value = -30486
$ws.Range("H88").Value = 144179.58
$ws.Range("I88").Value = 1239.2
$ws.Range("K88").Value = 1239.2
$ws.Range("M88").Value = -833.2
$ws.Range("H91").Value = 144179.58
$ws.Range("I91").Value = 1239.2
$ws.Range("K91").Value = 1239.2
$ws.Range("M91").Value = 164.8
$ws.Range("H132").Value = 8206739.5
$ws.Range("I132").Value = 10001633
$ws.Range("K132").Value = 30004899
$ws.Range("M132").Value = -30002369

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 962.5
$ws.Range("I16").Value = 847.2222
$ws.Range("J16").Value = 2000
$ws.Range("K16").Value = 847.2222
$ws.Range("L16").Value = 2000
$ws.Range("M16").Value = -560.2222
$ws.Range("N16").Value = -2574
$ws.Range("H58").Value = 20109.777
$ws.Range("I58").Value = 1794.5834
$ws.Range("J58").Value = 34761.934
$ws.Range("K58").Value = 1794.5834
$ws.Range("L58").Value = 34761.934
$ws.Range("M58").Value = -1591.5834
$ws.Range("N58").Value = -35167.934
$ws.Range("H86").Value = 13506.2
$ws.Range("I86").Value = 2499.5
$ws.Range("J86").Value = 16257.875
$ws.Range("K86").Value = 2499.5
$ws.Range("L86").Value = 16257.875
$ws.Range("M86").Value = -1376.5
$ws.Range("N86").Value = -18503.875
$ws.Range("H89").Value = 13506.2
$ws.Range("I89").Value = 2499.5
$ws.Range("J89").Value = 16257.875
$ws.Range("K89").Value = 12497.5
$ws.Range("L89").Value = 81289.375
$ws.Range("M89").Value = -6881.5
$ws.Range("N89").Value = -92521.375
$ws.Range("H113").Value = 962.5
$ws.Range("I113").Value = 847.2222
$ws.Range("J113").Value = 2000
$ws.Range("K113").Value = 847.2222
$ws.Range("L113").Value = 2000
$ws.Range("M113").Value = 1322.7778
$ws.Range("N113").Value = -6340
$ws.Range("H132").Value = 76927336
$ws.Range("I132").Value = 90911760
$ws.Range("J132").Value = 13006.5
$ws.Range("K132").Value = 272735280
$ws.Range("L132").Value = 39019.5
$ws.Range("M132").Value = -272732750
$ws.Range("N132").Value = -44079.5
$ws.Range("H136").Value = 20109.777
$ws.Range("I136").Value = 1794.5834
$ws.Range("J136").Value = 34761.934
$ws.Range("K136").Value = 5383.7502
$ws.Range("L136").Value = 104285.802
$ws.Range("M136").Value = -2833.7502
$ws.Range("N136").Value = -109385.802
$ws.Range("H140").Value = 50000
$ws.Range("J140").Value = 50000
$ws.Range("L140").Value = 50000
$ws.Range("N140").Value = -60360
$ws.Range("H141").Value = 18741.75
$ws.Range("J141").Value = 20704.857
$ws.Range("L141").Value = 20704.857
$ws.Range("N141").Value = -31064.857

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H82").Value = 10000
$ws.Range("I82").Value = 0
$ws.Range("K82").Value = 0
$ws.Range("M82").ClearContents()
$ws.Range("H85").Value = 10000
$ws.Range("I85").Value = 0
$ws.Range("K85").Value = 0
$ws.Range("M85").ClearContents()
$ws.Range("H103").Value = 1582.25
$ws.Range("I103").Value = 536
$ws.Range("J103").Value = 3326
$ws.Range("K103").Value = 1608
$ws.Range("L103").Value = 9978
$ws.Range("M103").Value = -729
$ws.Range("N103").Value = -11736
$ws.Range("H113").Value = 767.7
$ws.Range("I113").Value = 698.8333
$ws.Range("J113").Value = 797.2143
$ws.Range("K113").Value = 2096.4999
$ws.Range("L113").Value = 2391.6429
$ws.Range("M113").Value = 73.5001000000002
$ws.Range("N113").Value = -6731.6429
$ws.Range("H131").Value = 691.53
$ws.Range("I131").Value = 338
$ws.Range("J131").Value = 726.4945
$ws.Range("K131").Value = 1014
$ws.Range("L131").Value = 2179.4835
$ws.Range("M131").Value = 4026
$ws.Range("N131").Value = -12259.4835
$ws.Range("H132").Value = 991
$ws.Range("I132").Value = 998.6667
$ws.Range("K132").Value = 8988.0003
$ws.Range("M132").Value = -6458.0003

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 2608470.5
$ws.Range("I70").Value = 4357.0835
$ws.Range("J70").Value = 5212584
$ws.Range("K70").Value = 4357.0835
$ws.Range("L70").Value = 5212584
$ws.Range("M70").Value = -4087.0835
$ws.Range("N70").Value = -5213124
$ws.Range("H73").Value = 2608470.5
$ws.Range("I73").Value = 4357.0835
$ws.Range("J73").Value = 5212584
$ws.Range("K73").Value = 4357.0835
$ws.Range("L73").Value = 5212584
$ws.Range("M73").Value = -3421.0835
$ws.Range("N73").Value = -5214456
$ws.Range("H113").Value = 7004.3125
$ws.Range("I113").Value = 11897.375
$ws.Range("J113").Value = 2111.25
$ws.Range("K113").Value = 11897.375
$ws.Range("L113").Value = 2111.25
$ws.Range("M113").Value = -9727.375
$ws.Range("N113").Value = -6451.25
$ws.Range("H126").Value = 3369.5334
$ws.Range("J126").Value = 4589.5386
$ws.Range("L126").Value = 13768.6158
$ws.Range("N126").Value = -18708.6158

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H107").Value = 1750
$ws.Range("I107").Value = 1750
$ws.Range("K107").Value = 1750
$ws.Range("M107").Value = 170
$ws.Range("H136").Value = 1899.7222
$ws.Range("I136").Value = 1882.5714
$ws.Range("J136").Value = 2500
$ws.Range("K136").Value = 5647.7142
$ws.Range("L136").Value = 7500
$ws.Range("M136").Value = -3097.7142
$ws.Range("N136").Value = -12600

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 15625984
$ws.Range("I132").Value = 19231410
$ws.Range("K132").Value = 57694230
$ws.Range("M132").Value = -57691700
